# Regenerate / fix tests:
# 1. Insert a new "is_targeted list" sheet (TRUE/FALSE) right before "ms_source list".
# 2. Point the is_targeted (column N) data validation at the new list sheet instead of
#    the inline "TRUE,FALSE" literal list.

$wb = $excel.ActiveWorkbook

$mainSheet = $wb.Sheets.Item("Export as TSV")
$msSourceSheet = $wb.Sheets.Item("ms_source list")

# Insert the new lookup sheet directly before "ms_source list" so the tab order becomes:
# ... analyte_class list, is_targeted list, ms_source list, polarity list, ...
$isTargetedSheet = $wb.Worksheets.Add($msSourceSheet)
$isTargetedSheet.Name = "is_targeted list"

# Use a leading apostrophe so Excel stores these as plain text ("TRUE"/"FALSE") instead
# of native booleans, then clear the resulting quote-prefix formatting so no style is left
# behind on the cells.
$isTargetedSheet.Range("A1").Value = "'TRUE"
$isTargetedSheet.Range("A2").Value = "'FALSE"
$isTargetedSheet.Range("A1:A2").ClearFormats()

# Re-point column N's data validation at the new list sheet instead of the inline list.
$nRange = $mainSheet.Range("N2:N1048576")
$nRange.Validation.Delete()
$nRange.Validation.Add(3, 1, 1, "='is_targeted list'!`$A`$1:`$A`$2")
$nRange.Validation.ErrorTitle = "Value must come from list"
$nRange.Validation.ErrorMessage = "Value must be one of: TRUE / FALSE."
$nRange.Validation.IgnoreBlank = $true
$nRange.Validation.InCellDropdown = $true
$nRange.Validation.ShowInput = $true
$nRange.Validation.ShowError = $true

# Restore the originally active sheet/tab selection.
$mainSheet.Activate()
$mainSheet.Range("A1").Select()
